$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 174.22223
$ws.Range("I6").Value = 94.545456
$ws.Range("J6").Value = 299.42856
$ws.Range("K6").Value = 283.636368
$ws.Range("L6").Value = 898.28568
$ws.Range("M6").Value = -171.636368
$ws.Range("N6").Value = -1122.28568

$ws.Range("H40").Value = 44689344
$ws.Range("I40").Value = 6251595
$ws.Range("J40").Value = 72144880
$ws.Range("K40").Value = 6251595
$ws.Range("L40").Value = 72144880
$ws.Range("M40").Value = -6251420
$ws.Range("N40").Value = -72145230

$ws.Range("H74").Value = 3090.923
$ws.Range("I74").Value = 2743.818
$ws.Range("K74").Value = 2743.818
$ws.Range("M74").Value = -1807.818

$ws.Range("H77").Value = 3090.923
$ws.Range("I77").Value = 2743.818
$ws.Range("K77").Value = 13719.09
$ws.Range("M77").Value = -9039.09

$ws.Range("H116").Value = 2003.3334
$ws.Range("I116").Value = 2003.3334
$ws.Range("K116").Value = 2003.3334
$ws.Range("M116").Value = 1438.6666

$ws.Range("H132").Value = 11370488
$ws.Range("I132").Value = 13518930
$ws.Range("J132").Value = 14433.571
$ws.Range("K132").Value = 40556790
$ws.Range("L132").Value = 43300.713
$ws.Range("M132").Value = -40554260
$ws.Range("N132").Value = -48360.713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2310.4614
$ws.Range("I61").Value = 1748
$ws.Range("J61").Value = 2966.6667
$ws.Range("K61").Value = 1748
$ws.Range("L61").Value = 2966.6667
$ws.Range("M61").Value = -1536
$ws.Range("N61").Value = -3390.6667

$ws.Range("H74").Value = 71430216
$ws.Range("I74").Value = 83335000
$ws.Range("J74").Value = 1499
$ws.Range("K74").Value = 83335000
$ws.Range("L74").Value = 1499
$ws.Range("M74").Value = -83334126
$ws.Range("N74").Value = -3247

$ws.Range("H77").Value = 71430216
$ws.Range("I77").Value = 83335000
$ws.Range("J77").Value = 1499
$ws.Range("K77").Value = 416675000
$ws.Range("L77").Value = 7495
$ws.Range("M77").Value = -416670632
$ws.Range("N77").Value = -16231

$ws.Range("H136").Value = 2310.4614
$ws.Range("I136").Value = 1748
$ws.Range("J136").Value = 2966.6667
$ws.Range("K136").Value = 5244
$ws.Range("L136").Value = 8900.000100000001
$ws.Range("M136").Value = -2694
$ws.Range("N136").Value = -14000.0001

$ws.Range("H139").Value = 36019.285
$ws.Range("J139").Value = 36019.285
$ws.Range("L139").Value = 36019.285
$ws.Range("N139").Value = -46299.285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4635502.5
$ws.Range("I134").Value = 1327
$ws.Range("J134").Value = 8556728
$ws.Range("K134").Value = 3981
$ws.Range("L134").Value = 25670184
$ws.Range("M134").Value = -1446
$ws.Range("N134").Value = -25675254

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 90911470
$ws.Range("I99").Value = 200001760
$ws.Range("J99").Value = 2898.5
$ws.Range("K99").Value = 200001760
$ws.Range("L99").Value = 2898.5
$ws.Range("M99").Value = -200000262
$ws.Range("N99").Value = -5894.5

$ws.Range("H122").Value = 12500645
$ws.Range("I122").Value = 25000562
$ws.Range("J122").Value = 727.3
$ws.Range("K122").Value = 75001686
$ws.Range("L122").Value = 2181.9
$ws.Range("M122").Value = -74999236
$ws.Range("N122").Value = -7081.9

$ws.Range("H126").Value = 90911470
$ws.Range("I126").Value = 200001760
$ws.Range("J126").Value = 2898.5
$ws.Range("K126").Value = 600005280
$ws.Range("L126").Value = 8695.5
$ws.Range("M126").Value = -600002810
$ws.Range("N126").Value = -13635.5

$ws.Range("H132").Value = 11906542
$ws.Range("I132").Value = 1452.6818
$ws.Range("J132").Value = 55558536
$ws.Range("K132").Value = 4358.0454
$ws.Range("L132").Value = 166675608
$ws.Range("M132").Value = -1828.0454
$ws.Range("N132").Value = -166680668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 155
$ws.Range("J7").Value = 155
$ws.Range("L7").Value = 465
$ws.Range("N7").Value = -689

$ws.Range("H80").Value = 1257.8889
$ws.Range("J80").Value = 1257.8889
$ws.Range("L80").Value = 3773.6667
$ws.Range("N80").Value = -5645.6667

$ws.Range("H83").Value = 1257.8889
$ws.Range("J83").Value = 1257.8889
$ws.Range("L83").Value = 11321.0001
$ws.Range("N83").Value = -20681.0001

$ws.Range("H92").Value = 12364.889
$ws.Range("I92").Value = 300
$ws.Range("K92").Value = 900
$ws.Range("M92").Value = 348

$ws.Range("H104").Value = 37757.145
$ws.Range("I104").Value = 999
$ws.Range("J104").Value = 43883.5
$ws.Range("K104").Value = 2997
$ws.Range("L104").Value = 131650.5
$ws.Range("M104").Value = -376
$ws.Range("N104").Value = -136892.5

$ws.Range("H107").Value = 34486500
$ws.Range("I107").Value = 122.22222
$ws.Range("J107").Value = 50005372
$ws.Range("K107").Value = 366.66666
$ws.Range("L107").Value = 150016116
$ws.Range("M107").Value = 1553.33334
$ws.Range("N107").Value = -150019956

$ws.Range("H113").Value = 34444924
$ws.Range("I113").Value = 41667020
$ws.Range("J113").Value = 33333834
$ws.Range("K113").Value = 125001060
$ws.Range("L113").Value = 100001502
$ws.Range("M113").Value = -124998890
$ws.Range("N113").Value = -100005842

$ws.Range("H131").Value = 805.75
$ws.Range("I131").Value = 511.66666
$ws.Range("J131").Value = 824.5213
$ws.Range("K131").Value = 1534.99998
$ws.Range("L131").Value = 2473.5639
$ws.Range("M131").Value = 3505.00002
$ws.Range("N131").Value = -12553.5639

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 51021344
$ws.Range("I16").Value = 4202331
$ws.Range("J16").Value = 250002140
$ws.Range("K16").Value = 4202331
$ws.Range("L16").Value = 250002140
$ws.Range("M16").Value = -4202161
$ws.Range("N16").Value = -250002480

$ws.Range("H40").Value = 62500700
$ws.Range("I40").Value = 933.6667
$ws.Range("J40").Value = 250000000
$ws.Range("K40").Value = 933.6667
$ws.Range("L40").Value = 250000000
$ws.Range("M40").Value = -797.6667
$ws.Range("N40").Value = -250000272

$ws.Range("H82").Value = 1505.091
$ws.Range("I82").Value = 1272.5714
$ws.Range("J82").Value = 1912
$ws.Range("K82").Value = 1272.5714
$ws.Range("L82").Value = 1912
$ws.Range("M82").Value = -911.5714
$ws.Range("N82").Value = -2634

$ws.Range("H85").Value = 1505.091
$ws.Range("I85").Value = 1272.5714
$ws.Range("J85").Value = 1912
$ws.Range("K85").Value = 1272.5714
$ws.Range("L85").Value = 1912
$ws.Range("M85").Value = -24.57140000000004
$ws.Range("N85").Value = -4408

$ws.Range("H132").Value = 24397684
$ws.Range("I132").Value = 37039460
$ws.Range("J132").Value = 17114.5
$ws.Range("K132").Value = 111118380
$ws.Range("L132").Value = 51343.5
$ws.Range("M132").Value = -111115850
$ws.Range("N132").Value = -56403.5

$ws.Range("H136").Value = 81637336
$ws.Range("I136").Value = 10994058
$ws.Range("K136").Value = 32982174
$ws.Range("M136").Value = -32979624

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 36901.53
$ws.Range("I132").Value = 43885.043
$ws.Range("J132").Value = 15951
$ws.Range("K132").Value = 131655.129
$ws.Range("L132").Value = 47853
$ws.Range("M132").Value = -129125.129
$ws.Range("N132").Value = -52913

$ws.Range("H136").Value = 3247.6924
$ws.Range("I136").Value = 3950.4546
$ws.Range("J136").Value = 2027.1052
$ws.Range("K136").Value = 11851.3638
$ws.Range("L136").Value = 6081.3156
$ws.Range("M136").Value = -9301.363799999999
$ws.Range("N136").Value = -11181.3156
